# "More Clean up of Advanced Practices"
#
# On the "Dependency Injection" slide, the bullet that used to read
#   "Constructor Injection (preferred)"
# is reworded/split so that "(preferred)" becomes "(my preferred)", and the
# text ends up split across three runs:
#   "Constructor Injection " + "(my preferred" + ")"

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$shape = $s.Shapes.Item(2)
$tf = $shape.TextFrame
$tr = $tf.TextRange

$para1 = $tr.Paragraphs(1, 1)

# Replace "(preferred)" (chars 23-33 of "Constructor Injection (preferred)")
# with "(my preferred)" -- this naturally splits the paragraph into two runs:
# "Constructor Injection " and "(my preferred)".
$tail = $para1.Characters(23, 11)
$tail.Text = "(my preferred)"

# Now split the trailing ")" off into its own run by re-assigning just the
# last character back to itself.
$closingParen = $para1.Characters(36, 1)
$closingParen.Text = ")"
